$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), matching the new
# reporting-year column added to the EVOP yearly financials sheet.
$ws.Columns("D:D").Insert()

# The newly inserted column D has no number formatting; copy it over from
# column E (which holds the column that used to be D) so date/number
# formats line up with the rest of the table. Limit to the used range so we
# don't blow up the sheet with full-column copy/paste.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)

# Populate the new column D with the latest period (FY ending 2018-12-31,
# serial 43465) values for every data row.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 564800
$ws.Range("D9").Value = 189400
$ws.Range("D10").Value = 375400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 6200
$ws.Range("D15").Value = 87200
$ws.Range("D17").Value = 594100
$ws.Range("D18").Value = -29400
$ws.Range("D20").Value = 700
$ws.Range("D21").Value = 58500
$ws.Range("D22").Value = 59800
$ws.Range("D23").Value = -88400
$ws.Range("D24").Value = 10400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -98900
$ws.Range("D27").Value = -14700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -700
$ws.Range("D33").Value = -14700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -14700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 350700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 71600
$ws.Range("D44").Value = 8900
$ws.Range("D45").Value = 260100
$ws.Range("D46").Value = 691300
$ws.Range("D47").Value = 2700
$ws.Range("D48").Value = 103000
$ws.Range("D49").Value = 643100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 94200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1534400
$ws.Range("D57").Value = 48900
$ws.Range("D58").Value = 49000
$ws.Range("D59").Value = 545400
$ws.Range("D60").Value = 643400
$ws.Range("D61").Value = 676900
$ws.Range("D62").Value = 66700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1583000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -223800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -48600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -14700
$ws.Range("D83").Value = 87200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 202000
$ws.Range("D91").Value = -48800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -125600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 80600
$ws.Range("D101").Value = -11500
$ws.Range("D102").Value = 145600

# Row 91 (Capital Expenditures) received freshly-updated figures for the
# existing years as well, not just a column shift.
$ws.Range("E91").Value = -42000
$ws.Range("F91").Value = -31700
$ws.Range("G91").Value = -30200
